$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(32, 8).Value = 4500
$ws.Cells.Item(32, 10).Value = 4000
$ws.Cells.Item(32, 12).Value = 4000
$ws.Cells.Item(32, 14).Value = -4652

$ws.Cells.Item(33, 8).Value = 235.42308
$ws.Cells.Item(33, 9).Value = 238.41667
$ws.Cells.Item(33, 11).Value = 238.41667
$ws.Cells.Item(33, 13).Value = -9.416670000000011

$ws.Cells.Item(58, 8).Value = 5316.7144
$ws.Cells.Item(58, 10).Value = 6169.5
$ws.Cells.Item(58, 12).Value = 18508.5
$ws.Cells.Item(58, 14).Value = -18808.5

$ws.Cells.Item(62, 8).Value = 3230.6956
$ws.Cells.Item(62, 9).Value = 1991.3125
$ws.Cells.Item(62, 10).Value = 6063.5713
$ws.Cells.Item(62, 11).Value = 1991.3125
$ws.Cells.Item(62, 12).Value = 6063.5713
$ws.Cells.Item(62, 13).Value = -1367.3125
$ws.Cells.Item(62, 14).Value = -7311.5713

$ws.Cells.Item(65, 8).Value = 3230.6956
$ws.Cells.Item(65, 9).Value = 1991.3125
$ws.Cells.Item(65, 10).Value = 6063.5713
$ws.Cells.Item(65, 11).Value = 9956.5625
$ws.Cells.Item(65, 12).Value = 30317.8565
$ws.Cells.Item(65, 13).Value = -6836.5625
$ws.Cells.Item(65, 14).Value = -36557.85649999999

$ws.Cells.Item(69, 8).Value = 18731
$ws.Cells.Item(69, 9).Value = 6753.5
$ws.Cells.Item(69, 10).Value = 20908.727
$ws.Cells.Item(69, 11).Value = 20260.5
$ws.Cells.Item(69, 12).Value = 62726.181
$ws.Cells.Item(69, 13).Value = -19386.5
$ws.Cells.Item(69, 14).Value = -64474.181

$ws.Cells.Item(72, 8).Value = 18731
$ws.Cells.Item(72, 9).Value = 6753.5
$ws.Cells.Item(72, 10).Value = 20908.727
$ws.Cells.Item(72, 11).Value = 60781.5
$ws.Cells.Item(72, 12).Value = 188178.543
$ws.Cells.Item(72, 13).Value = -56413.5
$ws.Cells.Item(72, 14).Value = -196914.543

$ws.Cells.Item(88, 8).Value = 1514.909
$ws.Cells.Item(88, 10).Value = 962.6667
$ws.Cells.Item(88, 12).Value = 962.6667
$ws.Cells.Item(88, 14).Value = -1774.6667

$ws.Cells.Item(91, 8).Value = 1514.909
$ws.Cells.Item(91, 10).Value = 962.6667
$ws.Cells.Item(91, 12).Value = 962.6667
$ws.Cells.Item(91, 14).Value = -3770.6667

$ws.Cells.Item(97, 8).Value = 1333
$ws.Cells.Item(97, 10).Value = 1333
$ws.Cells.Item(97, 12).Value = 3999
$ws.Cells.Item(97, 14).Value = -4991

$ws.Cells.Item(106, 8).Value = 4899.381
$ws.Cells.Item(106, 9).Value = 5205.316
$ws.Cells.Item(106, 11).Value = 5205.316
$ws.Cells.Item(106, 13).Value = -4574.316

$ws.Cells.Item(113, 8).Value = 6264.5
$ws.Cells.Item(113, 9).Value = 6160.4
$ws.Cells.Item(113, 10).Value = 6368.6
$ws.Cells.Item(113, 11).Value = 6160.4
$ws.Cells.Item(113, 12).Value = 6368.6
$ws.Cells.Item(113, 13).Value = -2906.4
$ws.Cells.Item(113, 14).Value = -12876.6

$ws.Cells.Item(116, 8).Value = 13363.15
$ws.Cells.Item(116, 9).Value = 8440.666999999999
$ws.Cells.Item(116, 10).Value = 17390.637
$ws.Cells.Item(116, 11).Value = 8440.666999999999
$ws.Cells.Item(116, 12).Value = 17390.637
$ws.Cells.Item(116, 13).Value = -4998.666999999999
$ws.Cells.Item(116, 14).Value = -24274.637

$ws.Cells.Item(132, 8).Value = 2104.8667
$ws.Cells.Item(132, 9).Value = 1736.1282
$ws.Cells.Item(132, 11).Value = 5208.3846
$ws.Cells.Item(132, 13).Value = -2678.3846

$ws.Cells.Item(134, 8).Value = 89997
$ws.Cells.Item(134, 10).Value = 89997
$ws.Cells.Item(134, 12).Value = 89997
$ws.Cells.Item(134, 14).Value = -100137

$ws.Cells.Item(135, 8).Value = 909.7
$ws.Cells.Item(135, 9).Value = 897.9375
$ws.Cells.Item(135, 10).Value = 956.75
$ws.Cells.Item(135, 11).Value = 8081.4375
$ws.Cells.Item(135, 12).Value = 8610.75
$ws.Cells.Item(135, 13).Value = -5546.4375
$ws.Cells.Item(135, 14).Value = -13680.75

$ws.Cells.Item(137, 8).Value = 13515272
$ws.Cells.Item(137, 9).Value = 20001150
$ws.Cells.Item(137, 10).Value = 3029.1667
$ws.Cells.Item(137, 11).Value = 60003450
$ws.Cells.Item(137, 12).Value = 9087.500100000001
$ws.Cells.Item(137, 13).Value = -60000900
$ws.Cells.Item(137, 14).Value = -14187.5001

$ws.Cells.Item(138, 8).Value = 5778.8286
$ws.Cells.Item(138, 9).Value = 3512.125
$ws.Cells.Item(138, 10).Value = 6450.4443
$ws.Cells.Item(138, 11).Value = 10536.375
$ws.Cells.Item(138, 12).Value = 19351.3329
$ws.Cells.Item(138, 13).Value = -5396.375
$ws.Cells.Item(138, 14).Value = -29631.3329

$ws.Cells.Item(141, 8).Value = 22735592
$ws.Cells.Item(141, 9).Value = 26320194
$ws.Cells.Item(141, 10).Value = 33110.668
$ws.Cells.Item(141, 11).Value = 78960582
$ws.Cells.Item(141, 12).Value = 99332.00399999999
$ws.Cells.Item(141, 13).Value = -78955402
$ws.Cells.Item(141, 14).Value = -109692.004

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 918.2
$ws.Cells.Item(2, 9).Value = 1035.6364
$ws.Cells.Item(2, 10).Value = 774.6667
$ws.Cells.Item(2, 11).Value = 1035.6364
$ws.Cells.Item(2, 12).Value = 774.6667
$ws.Cells.Item(2, 13).Value = -922.6364000000001
$ws.Cells.Item(2, 14).Value = -1000.6667

$ws.Cells.Item(32, 8).Value = 5478.8306
$ws.Cells.Item(32, 10).Value = 6299.8335
$ws.Cells.Item(32, 12).Value = 6299.8335
$ws.Cells.Item(32, 14).Value = -6873.8335

$ws.Cells.Item(45, 8).Value = 6946.5454
$ws.Cells.Item(45, 9).Value = 2828.6667
$ws.Cells.Item(45, 11).Value = 2828.6667
$ws.Cells.Item(45, 13).Value = -2451.6667

$ws.Cells.Item(47, 8).Value = 89500
$ws.Cells.Item(47, 10).Value = 89500
$ws.Cells.Item(47, 12).Value = 89500
$ws.Cells.Item(47, 14).Value = -90950

$ws.Cells.Item(49, 8).Value = 89950
$ws.Cells.Item(49, 10).Value = 89950
$ws.Cells.Item(49, 12).Value = 89950
$ws.Cells.Item(49, 14).Value = -90470

$ws.Cells.Item(58, 8).Value = 82250
$ws.Cells.Item(58, 9).Value = 82250
$ws.Cells.Item(58, 10).Value = 0
$ws.Cells.Item(58, 11).Value = 82250
$ws.Cells.Item(58, 12).Value = 0
$ws.Cells.Item(58, 13).Value = -81820
$ws.Cells.Item(58, 14).ClearContents()

$ws.Cells.Item(61, 8).Value = 15108988
$ws.Cells.Item(61, 9).Value = 18921248
$ws.Cells.Item(61, 11).Value = 18921248
$ws.Cells.Item(61, 13).Value = -18921036

$ws.Cells.Item(63, 8).Value = 2706.9167
$ws.Cells.Item(63, 9).Value = 2848.4
$ws.Cells.Item(63, 10).Value = 1999.5
$ws.Cells.Item(63, 11).Value = 2848.4
$ws.Cells.Item(63, 12).Value = 1999.5
$ws.Cells.Item(63, 13).Value = -2162.4
$ws.Cells.Item(63, 14).Value = -3371.5

$ws.Cells.Item(66, 8).Value = 2706.9167
$ws.Cells.Item(66, 9).Value = 2848.4
$ws.Cells.Item(66, 10).Value = 1999.5
$ws.Cells.Item(66, 11).Value = 14242
$ws.Cells.Item(66, 12).Value = 9997.5
$ws.Cells.Item(66, 13).Value = -10810
$ws.Cells.Item(66, 14).Value = -16861.5

$ws.Cells.Item(74, 8).Value = 1952.579
$ws.Cells.Item(74, 9).Value = 1810.7273
$ws.Cells.Item(74, 10).Value = 2147.625
$ws.Cells.Item(74, 11).Value = 1810.7273
$ws.Cells.Item(74, 12).Value = 2147.625
$ws.Cells.Item(74, 13).Value = -936.7273
$ws.Cells.Item(74, 14).Value = -3895.625

$ws.Cells.Item(77, 8).Value = 1952.579
$ws.Cells.Item(77, 9).Value = 1810.7273
$ws.Cells.Item(77, 10).Value = 2147.625
$ws.Cells.Item(77, 11).Value = 9053.636500000001
$ws.Cells.Item(77, 12).Value = 10738.125
$ws.Cells.Item(77, 13).Value = -4685.636500000001
$ws.Cells.Item(77, 14).Value = -19474.125

$ws.Cells.Item(88, 8).Value = 5053.4287
$ws.Cells.Item(88, 9).Value = 3449.5
$ws.Cells.Item(88, 10).Value = 5695
$ws.Cells.Item(88, 11).Value = 3449.5
$ws.Cells.Item(88, 12).Value = 5695
$ws.Cells.Item(88, 13).Value = -3043.5
$ws.Cells.Item(88, 14).Value = -6507

$ws.Cells.Item(91, 8).Value = 5053.4287
$ws.Cells.Item(91, 9).Value = 3449.5
$ws.Cells.Item(91, 10).Value = 5695
$ws.Cells.Item(91, 11).Value = 3449.5
$ws.Cells.Item(91, 12).Value = 5695
$ws.Cells.Item(91, 13).Value = -2045.5
$ws.Cells.Item(91, 14).Value = -8503

$ws.Cells.Item(97, 8).Value = 1771.6207
$ws.Cells.Item(97, 9).Value = 1780.6666
$ws.Cells.Item(97, 10).Value = 1756.8182
$ws.Cells.Item(97, 11).Value = 1780.6666
$ws.Cells.Item(97, 12).Value = 1756.8182
$ws.Cells.Item(97, 13).Value = -1284.6666
$ws.Cells.Item(97, 14).Value = -2748.8182

$ws.Cells.Item(107, 8).Value = 150000
$ws.Cells.Item(107, 10).Value = 150000
$ws.Cells.Item(107, 12).Value = 150000
$ws.Cells.Item(107, 14).Value = -157680

$ws.Cells.Item(116, 8).Value = 918.2
$ws.Cells.Item(116, 9).Value = 1035.6364
$ws.Cells.Item(116, 10).Value = 774.6667
$ws.Cells.Item(116, 11).Value = 1035.6364
$ws.Cells.Item(116, 12).Value = 774.6667
$ws.Cells.Item(116, 13).Value = 1258.3636
$ws.Cells.Item(116, 14).Value = -5362.6667

$ws.Cells.Item(132, 8).Value = 2382210.2
$ws.Cells.Item(132, 9).Value = 1288.5122
$ws.Cells.Item(132, 10).Value = 100000000
$ws.Cells.Item(132, 11).Value = 3865.536599999999
$ws.Cells.Item(132, 12).Value = 300000000
$ws.Cells.Item(132, 13).Value = -1335.536599999999
$ws.Cells.Item(132, 14).Value = -300005060

$ws.Cells.Item(136, 8).Value = 15108988
$ws.Cells.Item(136, 9).Value = 18921248
$ws.Cells.Item(136, 11).Value = 56763744
$ws.Cells.Item(136, 13).Value = -56761194

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 918.2
$ws.Cells.Item(3, 9).Value = 1035.6364
$ws.Cells.Item(3, 10).Value = 774.6667
$ws.Cells.Item(3, 11).Value = 1035.6364
$ws.Cells.Item(3, 12).Value = 774.6667
$ws.Cells.Item(3, 13).Value = -921.6364000000001
$ws.Cells.Item(3, 14).Value = -1002.6667

$ws.Cells.Item(20, 8).Value = 7736.7334
$ws.Cells.Item(20, 10).Value = 3859.75
$ws.Cells.Item(20, 12).Value = 3859.75
$ws.Cells.Item(20, 14).Value = -4353.75

$ws.Cells.Item(22, 8).Value = 2593.8333
$ws.Cells.Item(22, 9).Value = 1512.6
$ws.Cells.Item(22, 10).Value = 8000
$ws.Cells.Item(22, 11).Value = 1512.6
$ws.Cells.Item(22, 12).Value = 8000
$ws.Cells.Item(22, 13).Value = -1339.6
$ws.Cells.Item(22, 14).Value = -8346

$ws.Cells.Item(105, 8).Value = 604686.75
$ws.Cells.Item(105, 9).Value = 848948.25
$ws.Cells.Item(105, 10).Value = 5135.909
$ws.Cells.Item(105, 11).Value = 848948.25
$ws.Cells.Item(105, 12).Value = 5135.909
$ws.Cells.Item(105, 13).Value = -847201.25
$ws.Cells.Item(105, 14).Value = -8629.909

$ws.Cells.Item(113, 8).Value = 15040
$ws.Cells.Item(113, 9).Value = 15040
$ws.Cells.Item(113, 11).Value = 15040
$ws.Cells.Item(113, 13).Value = -12870

$ws.Cells.Item(134, 8).Value = 4548010.5
$ws.Cells.Item(134, 9).Value = 2123.0667
$ws.Cells.Item(134, 11).Value = 6369.2001
$ws.Cells.Item(134, 13).Value = -3834.2001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(16, 8).Value = 14286119
$ws.Cells.Item(16, 9).Value = 14286119
$ws.Cells.Item(16, 11).Value = 14286119
$ws.Cells.Item(16, 13).Value = -14285832

$ws.Cells.Item(31, 8).Value = 12823237
$ws.Cells.Item(31, 9).Value = 19233158
$ws.Cells.Item(31, 10).Value = 3395.7693
$ws.Cells.Item(31, 11).Value = 19233158
$ws.Cells.Item(31, 12).Value = 3395.7693
$ws.Cells.Item(31, 13).Value = -19232863
$ws.Cells.Item(31, 14).Value = -3985.7693

$ws.Cells.Item(34, 8).Value = 12823237
$ws.Cells.Item(34, 9).Value = 19233158
$ws.Cells.Item(34, 10).Value = 3395.7693
$ws.Cells.Item(34, 11).Value = 19233158
$ws.Cells.Item(34, 12).Value = 3395.7693
$ws.Cells.Item(34, 13).Value = -19232956
$ws.Cells.Item(34, 14).Value = -3799.7693

$ws.Cells.Item(86, 8).Value = 8749.210999999999
$ws.Cells.Item(86, 9).Value = 9725.532999999999
$ws.Cells.Item(86, 11).Value = 9725.532999999999
$ws.Cells.Item(86, 13).Value = -8602.532999999999

$ws.Cells.Item(89, 8).Value = 8749.210999999999
$ws.Cells.Item(89, 9).Value = 9725.532999999999
$ws.Cells.Item(89, 11).Value = 48627.66499999999
$ws.Cells.Item(89, 13).Value = -43011.66499999999

$ws.Cells.Item(105, 8).Value = 1540.8
$ws.Cells.Item(105, 9).Value = 1001.25
$ws.Cells.Item(105, 11).Value = 1001.25
$ws.Cells.Item(105, 13).Value = 745.75

$ws.Cells.Item(107, 8).Value = 785.8
$ws.Cells.Item(107, 9).Value = 500.65625
$ws.Cells.Item(107, 11).Value = 500.65625
$ws.Cells.Item(107, 13).Value = 1419.34375

$ws.Cells.Item(113, 8).Value = 14286119
$ws.Cells.Item(113, 9).Value = 14286119
$ws.Cells.Item(113, 11).Value = 14286119
$ws.Cells.Item(113, 13).Value = -14283949

$ws.Cells.Item(132, 8).Value = 2527.0908
$ws.Cells.Item(132, 9).Value = 1572.2222
$ws.Cells.Item(132, 11).Value = 4716.6666
$ws.Cells.Item(132, 13).Value = -2186.6666

$ws.Cells.Item(133, 8).Value = 64750
$ws.Cells.Item(133, 10).Value = 64750
$ws.Cells.Item(133, 12).Value = 64750
$ws.Cells.Item(133, 14).Value = -69810

$ws.Cells.Item(134, 8).Value = 1607.3823
$ws.Cells.Item(134, 9).Value = 1663.3549
$ws.Cells.Item(134, 10).Value = 1029
$ws.Cells.Item(134, 11).Value = 4990.0647
$ws.Cells.Item(134, 12).Value = 3087
$ws.Cells.Item(134, 13).Value = -2455.0647
$ws.Cells.Item(134, 14).Value = -8157

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 2389641.2
$ws.Cells.Item(33, 9).Value = 110.47369
$ws.Cells.Item(33, 10).Value = 6929749.5
$ws.Cells.Item(33, 11).Value = 662.84214
$ws.Cells.Item(33, 12).Value = 41578497
$ws.Cells.Item(33, 13).Value = -379.84214
$ws.Cells.Item(33, 14).Value = -41579063

$ws.Cells.Item(62, 8).Value = 14834
$ws.Cells.Item(62, 9).Value = 8011.5
$ws.Cells.Item(62, 11).Value = 24034.5
$ws.Cells.Item(62, 13).Value = -23348.5

$ws.Cells.Item(63, 8).Value = 25444.584
$ws.Cells.Item(63, 9).Value = 16669.334
$ws.Cells.Item(63, 11).Value = 50008.00199999999
$ws.Cells.Item(63, 13).Value = -49259.00199999999

$ws.Cells.Item(64, 8).Value = 10606
$ws.Cells.Item(64, 9).Value = 5681.8
$ws.Cells.Item(64, 10).Value = 22916.5
$ws.Cells.Item(64, 11).Value = 17045.4
$ws.Cells.Item(64, 12).Value = 68749.5
$ws.Cells.Item(64, 13).Value = -16775.4
$ws.Cells.Item(64, 14).Value = -69289.5

$ws.Cells.Item(65, 8).Value = 14834
$ws.Cells.Item(65, 9).Value = 8011.5
$ws.Cells.Item(65, 11).Value = 72103.5
$ws.Cells.Item(65, 13).Value = -68671.5

$ws.Cells.Item(66, 8).Value = 25444.584
$ws.Cells.Item(66, 9).Value = 16669.334
$ws.Cells.Item(66, 11).Value = 150024.006
$ws.Cells.Item(66, 13).Value = -146280.006

$ws.Cells.Item(67, 8).Value = 10606
$ws.Cells.Item(67, 9).Value = 5681.8
$ws.Cells.Item(67, 10).Value = 22916.5
$ws.Cells.Item(67, 11).Value = 17045.4
$ws.Cells.Item(67, 12).Value = 68749.5
$ws.Cells.Item(67, 13).Value = -16109.4
$ws.Cells.Item(67, 14).Value = -70621.5

$ws.Cells.Item(70, 8).Value = 16385.273
$ws.Cells.Item(70, 9).Value = 9381.200000000001
$ws.Cells.Item(70, 10).Value = 22222
$ws.Cells.Item(70, 11).Value = 28143.6
$ws.Cells.Item(70, 12).Value = 66666
$ws.Cells.Item(70, 13).Value = -27828.6
$ws.Cells.Item(70, 14).Value = -67296

$ws.Cells.Item(73, 8).Value = 16385.273
$ws.Cells.Item(73, 9).Value = 9381.200000000001
$ws.Cells.Item(73, 10).Value = 22222
$ws.Cells.Item(73, 11).Value = 28143.6
$ws.Cells.Item(73, 12).Value = 66666
$ws.Cells.Item(73, 13).Value = -27051.6
$ws.Cells.Item(73, 14).Value = -68850

$ws.Cells.Item(74, 8).Value = 27777.666
$ws.Cells.Item(74, 10).Value = 27777.666
$ws.Cells.Item(74, 12).Value = 83332.99800000001
$ws.Cells.Item(74, 14).Value = -85454.99800000001

$ws.Cells.Item(77, 8).Value = 27777.666
$ws.Cells.Item(77, 10).Value = 27777.666
$ws.Cells.Item(77, 12).Value = 249998.994
$ws.Cells.Item(77, 14).Value = -260606.994

$ws.Cells.Item(87, 8).Value = 12058.25
$ws.Cells.Item(87, 9).Value = 4966.6665
$ws.Cells.Item(87, 11).Value = 14899.9995
$ws.Cells.Item(87, 13).Value = -13651.9995

$ws.Cells.Item(90, 8).Value = 12058.25
$ws.Cells.Item(90, 9).Value = 4966.6665
$ws.Cells.Item(90, 11).Value = 44699.9985
$ws.Cells.Item(90, 13).Value = -38459.9985

$ws.Cells.Item(92, 8).Value = 173.75
$ws.Cells.Item(92, 9).Value = 100
$ws.Cells.Item(92, 10).Value = 198.33333
$ws.Cells.Item(92, 11).Value = 300
$ws.Cells.Item(92, 12).Value = 594.99999
$ws.Cells.Item(92, 13).Value = 948
$ws.Cells.Item(92, 14).Value = -3090.99999

$ws.Cells.Item(131, 8).Value = 4284.48
$ws.Cells.Item(131, 9).Value = 3269.7778
$ws.Cells.Item(131, 10).Value = 4855.25
$ws.Cells.Item(131, 11).Value = 9809.3334
$ws.Cells.Item(131, 12).Value = 14565.75
$ws.Cells.Item(131, 13).Value = -4769.3334
$ws.Cells.Item(131, 14).Value = -24645.75

$ws.Cells.Item(141, 8).Value = 11465.6
$ws.Cells.Item(141, 9).Value = 5998.75
$ws.Cells.Item(141, 11).Value = 17996.25
$ws.Cells.Item(141, 13).Value = -12816.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 10340.8
$ws.Cells.Item(70, 9).Value = 8546
$ws.Cells.Item(70, 10).Value = 11911.25
$ws.Cells.Item(70, 11).Value = 8546
$ws.Cells.Item(70, 12).Value = 11911.25
$ws.Cells.Item(70, 13).Value = -8276
$ws.Cells.Item(70, 14).Value = -12451.25

$ws.Cells.Item(73, 8).Value = 10340.8
$ws.Cells.Item(73, 9).Value = 8546
$ws.Cells.Item(73, 10).Value = 11911.25
$ws.Cells.Item(73, 11).Value = 8546
$ws.Cells.Item(73, 12).Value = 11911.25
$ws.Cells.Item(73, 13).Value = -7610
$ws.Cells.Item(73, 14).Value = -13783.25

$ws.Cells.Item(74, 8).Value = 99999
$ws.Cells.Item(74, 10).Value = 99999
$ws.Cells.Item(74, 12).Value = 99999
$ws.Cells.Item(74, 14).Value = -101871

$ws.Cells.Item(77, 8).Value = 99999
$ws.Cells.Item(77, 10).Value = 99999
$ws.Cells.Item(77, 12).Value = 299997
$ws.Cells.Item(77, 14).Value = -309357

$ws.Cells.Item(80, 8).Value = 24002204
$ws.Cells.Item(80, 9).Value = 1934.75
$ws.Cells.Item(80, 10).Value = 120003280
$ws.Cells.Item(80, 11).Value = 1934.75
$ws.Cells.Item(80, 12).Value = 120003280
$ws.Cells.Item(80, 13).Value = -936.75
$ws.Cells.Item(80, 14).Value = -120005276

$ws.Cells.Item(83, 8).Value = 24002204
$ws.Cells.Item(83, 9).Value = 1934.75
$ws.Cells.Item(83, 10).Value = 120003280
$ws.Cells.Item(83, 11).Value = 9673.75
$ws.Cells.Item(83, 12).Value = 600016400
$ws.Cells.Item(83, 13).Value = -4681.75
$ws.Cells.Item(83, 14).Value = -600026384

$ws.Cells.Item(102, 8).Value = 1973.3334
$ws.Cells.Item(102, 9).Value = 1667.8
$ws.Cells.Item(102, 11).Value = 1667.8
$ws.Cells.Item(102, 13).Value = -45.79999999999995

$ws.Cells.Item(132, 8).Value = 8290590
$ws.Cells.Item(132, 9).Value = 2426.75
$ws.Cells.Item(132, 10).Value = 46968690
$ws.Cells.Item(132, 11).Value = 7280.25
$ws.Cells.Item(132, 12).Value = 140906070
$ws.Cells.Item(132, 13).Value = -4750.25
$ws.Cells.Item(132, 14).Value = -140911130

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 2544.8386
$ws.Cells.Item(16, 9).Value = 1243.84
$ws.Cells.Item(16, 10).Value = 7965.6665
$ws.Cells.Item(16, 11).Value = 1243.84
$ws.Cells.Item(16, 12).Value = 7965.6665
$ws.Cells.Item(16, 13).Value = -1073.84
$ws.Cells.Item(16, 14).Value = -8305.666499999999

$ws.Cells.Item(22, 8).Value = 10338.23
$ws.Cells.Item(22, 9).Value = 15449.875
$ws.Cells.Item(22, 10).Value = 2159.6
$ws.Cells.Item(22, 11).Value = 15449.875
$ws.Cells.Item(22, 12).Value = 2159.6
$ws.Cells.Item(22, 13).Value = -15154.875
$ws.Cells.Item(22, 14).Value = -2749.6

$ws.Cells.Item(27, 8).Value = 10338.23
$ws.Cells.Item(27, 9).Value = 15449.875
$ws.Cells.Item(27, 10).Value = 2159.6
$ws.Cells.Item(27, 11).Value = 15449.875
$ws.Cells.Item(27, 12).Value = 2159.6
$ws.Cells.Item(27, 13).Value = -15342.875
$ws.Cells.Item(27, 14).Value = -2373.6

$ws.Cells.Item(61, 8).Value = 58827556
$ws.Cells.Item(61, 9).Value = 90912120
$ws.Cells.Item(61, 11).Value = 90912120
$ws.Cells.Item(61, 13).Value = -90911918

$ws.Cells.Item(113, 8).Value = 58827556
$ws.Cells.Item(113, 9).Value = 90912120
$ws.Cells.Item(113, 11).Value = 90912120
$ws.Cells.Item(113, 13).Value = -90909950

$ws.Cells.Item(122, 8).Value = 3840.465
$ws.Cells.Item(122, 9).Value = 3420.639
$ws.Cells.Item(122, 11).Value = 10261.917
$ws.Cells.Item(122, 13).Value = -7811.917000000001

$ws.Cells.Item(132, 8).Value = 3038.9038
$ws.Cells.Item(132, 9).Value = 1827.4474
$ws.Cells.Item(132, 10).Value = 6327.143
$ws.Cells.Item(132, 11).Value = 5482.3422
$ws.Cells.Item(132, 12).Value = 18981.429
$ws.Cells.Item(132, 13).Value = -2952.3422
$ws.Cells.Item(132, 14).Value = -24041.429

$ws.Cells.Item(136, 8).Value = 2734.4167
$ws.Cells.Item(136, 9).Value = 1087.8572
$ws.Cells.Item(136, 10).Value = 5039.6
$ws.Cells.Item(136, 11).Value = 3263.5716
$ws.Cells.Item(136, 12).Value = 15118.8
$ws.Cells.Item(136, 13).Value = -713.5715999999998
$ws.Cells.Item(136, 14).Value = -20218.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(18, 8).Value = 19013.766
$ws.Cells.Item(18, 9).Value = 18155
$ws.Cells.Item(18, 10).Value = 19979.875
$ws.Cells.Item(18, 11).Value = 18155
$ws.Cells.Item(18, 12).Value = 19979.875
$ws.Cells.Item(18, 13).Value = -17982
$ws.Cells.Item(18, 14).Value = -20325.875

$ws.Cells.Item(81, 8).Value = 2561.9285
$ws.Cells.Item(81, 9).Value = 2705.5386
$ws.Cells.Item(81, 11).Value = 5411.0772
$ws.Cells.Item(81, 13).Value = -4350.0772

$ws.Cells.Item(84, 8).Value = 2561.9285
$ws.Cells.Item(84, 9).Value = 2705.5386
$ws.Cells.Item(84, 11).Value = 27055.386
$ws.Cells.Item(84, 13).Value = -21751.386

$ws.Cells.Item(92, 8).Value = 65249.75
$ws.Cells.Item(92, 10).Value = 65249.75
$ws.Cells.Item(92, 12).Value = 65249.75
$ws.Cells.Item(92, 14).Value = -70241.75

$ws.Cells.Item(100, 8).Value = 1628.5
$ws.Cells.Item(100, 9).Value = 1993.125
$ws.Cells.Item(100, 11).Value = 3986.25
$ws.Cells.Item(100, 13).Value = -3445.25

$ws.Cells.Item(107, 8).Value = 3360.7112
$ws.Cells.Item(107, 9).Value = 2063.3914
$ws.Cells.Item(107, 11).Value = 6190.174199999999
$ws.Cells.Item(107, 13).Value = -4270.174199999999

$ws.Cells.Item(113, 8).Value = 608.8570999999999
$ws.Cells.Item(113, 9).Value = 574.9091
$ws.Cells.Item(113, 11).Value = 1724.7273
$ws.Cells.Item(113, 13).Value = 445.2727

$ws.Cells.Item(122, 8).Value = 3000.6667
$ws.Cells.Item(122, 9).Value = 3000.6667
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 9002.000100000001
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -6552.000100000001
$ws.Cells.Item(122, 14).ClearContents()

$ws.Cells.Item(126, 8).Value = 2980.5
$ws.Cells.Item(126, 9).Value = 3401.5715
$ws.Cells.Item(126, 11).Value = 10204.7145
$ws.Cells.Item(126, 13).Value = -7734.7145

$ws.Cells.Item(132, 8).Value = 1668598.4
$ws.Cells.Item(132, 9).Value = 2318
$ws.Cells.Item(132, 10).Value = 10000000
$ws.Cells.Item(132, 11).Value = 6954
$ws.Cells.Item(132, 12).Value = 30000000
$ws.Cells.Item(132, 13).Value = -4424
$ws.Cells.Item(132, 14).Value = -30005060

$ws.Cells.Item(141, 8).Value = 92991.8
$ws.Cells.Item(141, 10).Value = 92991.8
$ws.Cells.Item(141, 12).Value = 92991.8
$ws.Cells.Item(141, 14).Value = -103351.8
